$p = $ppt.ActivePresentation
$tm = $p.TitleMaster
Write-Output $tm.Name
